$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-28 09:41:58"
$wsZh.Range("G5").Value = "2016-01-28 09:42:44"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-28 09:42:11"
$wsDe.Range("G5").Value = "2016-01-28 09:43:08"
